$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027510331597459
$ws.Range("D2").Value = 1.028816560393554
$ws.Range("E2").Value = 1.036184467811835
$ws.Range("F2").Value = 1.043970251305802
$ws.Range("I2").Value = 1.023594999628091
$ws.Range("J2").Value = 1.032667803777056
$ws.Range("K2").Value = 1.031632276888562
$ws.Range("L2").Value = 1.038978935124798
$ws.Range("M2").Value = 1.046742618042306
$ws.Range("N2").Value = 1.034134310507094
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028893043923549
$ws.Range("D3").Value = 1.030039968870523
$ws.Range("E3").Value = 1.037455394122874
$ws.Range("F3").Value = 1.045422046661947
$ws.Range("I3").Value = 1.023504579208683
$ws.Range("J3").Value = 1.033688507682869
$ws.Range("K3").Value = 1.032662580953607
$ws.Range("L3").Value = 1.040058186330808
$ws.Range("M3").Value = 1.048003884081036
$ws.Range("N3").Value = 1.035156463929508
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029787481076281
$ws.Range("D4").Value = 1.030831580599395
$ws.Range("E4").Value = 1.038277790879478
$ws.Range("F4").Value = 1.046361769813064
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.034348301310005
$ws.Range("K4").Value = 1.033328686352143
$ws.Range("L4").Value = 1.040756008892513
$ws.Range("M4").Value = 1.048819810016611
$ws.Range("N4").Value = 1.035817194539274
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030163442453253
$ws.Range("D5").Value = 1.031164373793181
$ws.Range("E5").Value = 1.038623535380957
$ws.Range("F5").Value = 1.046756909016546
$ws.Range("I5").Value = 1.023416932628374
$ws.Range("J5").Value = 1.034625521237688
$ws.Range("K5").Value = 1.033608583556463
$ws.Range("L5").Value = 1.041049250978025
$ws.Range("M5").Value = 1.049162781780916
$ws.Range("N5").Value = 1.036094808151051
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030226564640103
$ws.Range("D6").Value = 1.031220251308486
$ws.Range("E6").Value = 1.038681588028801
$ws.Range("F6").Value = 1.046823259419736
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.034672058548404
$ws.Range("K6").Value = 1.03365557178718
$ws.Range("L6").Value = 1.041098480522787
$ws.Range("M6").Value = 1.049220365808535
$ws.Range("N6").Value = 1.036141411550087
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029792504920548
$ws.Range("D7").Value = 1.030836027395722
$ws.Range("E7").Value = 1.038282410695053
$ws.Range("F7").Value = 1.04636704936435
$ws.Range("I7").Value = 1.023442985510381
$ws.Range("J7").Value = 1.03435200615054
$ws.Range("K7").Value = 1.033332426873142
$ws.Range("L7").Value = 1.040759927685062
$ws.Range("M7").Value = 1.048824392988327
$ws.Range("N7").Value = 1.035820904641107
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027977684463045
$ws.Range("D8").Value = 1.029230022185553
$ws.Range("E8").Value = 1.036613980787313
$ws.Range("F8").Value = 1.044460829820963
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.033012896046095
$ws.Range("K8").Value = 1.031980592852108
$ws.Range("L8").Value = 1.039343783782488
$ws.Range("M8").Value = 1.047168911987313
$ws.Range("N8").Value = 1.034479892846716
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02477746014151
$ws.Range("D9").Value = 1.026399770131655
$ws.Range("E9").Value = 1.03367400778006
$ws.Range("F9").Value = 1.041104043147697
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.030647942229689
$ws.Range("K9").Value = 1.029593986455951
$ws.Range("L9").Value = 1.036844191959806
$ws.Range("M9").Value = 1.044250060285979
$ws.Range("N9").Value = 1.032111580524626
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.022642151200193
$ws.Range("D10").Value = 1.024512548296931
$ws.Range("E10").Value = 1.031713812199101
$ws.Range("F10").Value = 1.038867404259756
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.029067567110577
$ws.Range("K10").Value = 1.02799969860244
$ws.Range("L10").Value = 1.0351747950427
$ws.Range("M10").Value = 1.042302789122537
$ws.Range("N10").Value = 1.030528961091574
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021717040028187
$ws.Range("D11").Value = 1.023695221703368
$ws.Range("E11").Value = 1.030864925184076
$ws.Range("F11").Value = 1.037899143375256
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.028382316426071
$ws.Range("K11").Value = 1.027308548652398
$ws.Range("L11").Value = 1.034451170293381
$ws.Range("M11").Value = 1.041459220739038
$ws.Range("N11").Value = 1.029842737272496
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.021373330462032
$ws.Range("D12").Value = 1.023391604010615
$ws.Range("E12").Value = 1.0305495899387
$ws.Range("F12").Value = 1.037539515714689
$ws.Range("I12").Value = 1.023938164268129
$ws.Range("J12").Value = 1.028127639064291
$ws.Range("K12").Value = 1.027051698930342
$ws.Range("L12").Value = 1.034182264958065
$ws.Range("M12").Value = 1.041145819046869
$ws.Range("N12").Value = 1.029587698239651
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.021447061137352
$ws.Range("D13").Value = 1.023456732288299
$ws.Range("E13").Value = 1.030617231431484
$ws.Range("F13").Value = 1.037616655907344
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.028182274834975
$ws.Range("K13").Value = 1.027106799823218
$ws.Range("L13").Value = 1.034239951509002
$ws.Range("M13").Value = 1.041213047692358
$ws.Range("N13").Value = 1.029642411599396
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021688630595918
$ws.Range("D14").Value = 1.023670125128551
$ws.Range("E14").Value = 1.030838859924248
$ws.Range("F14").Value = 1.03786941590941
$ws.Range("I14").Value = 1.023922712353276
$ws.Range("J14").Value = 1.028361267691008
$ws.Range("K14").Value = 1.027287319963414
$ws.Range("L14").Value = 1.034428944942706
$ws.Range("M14").Value = 1.041433316166309
$ws.Range("N14").Value = 1.029821658645816
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02183745846347
$ws.Range("D15").Value = 1.02380159990479
$ws.Range("E15").Value = 1.030975409664376
$ws.Range("F15").Value = 1.038025153180162
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.028471531801368
$ws.Range("K15").Value = 1.027398527590114
$ws.Range("L15").Value = 1.034545374162416
$ws.Range("M15").Value = 1.041569022376538
$ws.Range("N15").Value = 1.029932079343859
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.022703536584787
$ws.Range("D16").Value = 1.024566788122119
$ws.Range("E16").Value = 1.031770147356644
$ws.Range("F16").Value = 1.038931668698935
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.029113024854575
$ws.Range("K16").Value = 1.028045550504569
$ws.Range("L16").Value = 1.035222803180625
$ws.Range("M16").Value = 1.042358765423203
$ws.Range("N16").Value = 1.030574483390782
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023246664388997
$ws.Range("D17").Value = 1.025046728201883
$ws.Range("E17").Value = 1.032268632942378
$ws.Range("F17").Value = 1.039500356710896
$ws.Range("I17").Value = 1.023843013862359
$ws.Range("J17").Value = 1.029515162702889
$ws.Range("K17").Value = 1.028451190638781
$ws.Range("L17").Value = 1.035647528930459
$ws.Range("M17").Value = 1.042854043434038
$ws.Range("N17").Value = 1.030977192320967
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023563412991789
$ws.Range("D18").Value = 1.025326655217996
$ws.Range("E18").Value = 1.032559380893347
$ws.Range("F18").Value = 1.039832084249132
$ws.Range("I18").Value = 1.023826118601784
$ws.Range("J18").Value = 1.029749632680743
$ws.Range("K18").Value = 1.028687715498334
$ws.Range("L18").Value = 1.035895190745699
$ws.Range("M18").Value = 1.043142893684338
$ws.Range("N18").Value = 1.031211995273084
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023671407999188
$ws.Range("D19").Value = 1.02542210086033
$ws.Range("E19").Value = 1.032658516822929
$ws.Range("F19").Value = 1.039945198563256
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.029829565696023
$ws.Range("K19").Value = 1.028768351217202
$ws.Range("L19").Value = 1.035979624672343
$ws.Range("M19").Value = 1.043241378014305
$ws.Range("N19").Value = 1.031292041802413
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023188396991427
$ws.Range("D20").Value = 1.02499523664348
$ws.Range("E20").Value = 1.032215151250651
$ws.Range("F20").Value = 1.039439339691132
$ws.Range("I20").Value = 1.023846096030192
$ws.Range("J20").Value = 1.029472026479032
$ws.Range("K20").Value = 1.02840767740486
$ws.Range("L20").Value = 1.035601967512463
$ws.Range("M20").Value = 1.042800908656557
$ws.Range("N20").Value = 1.030933994838725
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021617496723348
$ws.Range("D21").Value = 1.023607287012213
$ws.Range("E21").Value = 1.030773596466804
$ws.Range("F21").Value = 1.037794983670888
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.028308562782892
$ws.Range("K21").Value = 1.0272341647859
$ws.Range("L21").Value = 1.034373294400773
$ws.Range("M21").Value = 1.041368454362953
$ws.Range("N21").Value = 1.029768878890684
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020629331602537
$ws.Range("D22").Value = 1.022734474334227
$ws.Range("E22").Value = 1.029867111170914
$ws.Range("F22").Value = 1.036761268503874
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 1.027576207601987
$ws.Range("K22").Value = 1.026495600789229
$ws.Range("L22").Value = 1.033600089044928
$ws.Range("M22").Value = 1.040467447624599
$ws.Range("N22").Value = 1.029035483681421
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.021153223831285
$ws.Range("D23").Value = 1.023197184788923
$ws.Range("E23").Value = 1.03034766926263
$ws.Range("F23").Value = 1.037309247385818
$ws.Range("I23").Value = 1.023948818288664
$ws.Range("J23").Value = 1.02796452379171
$ws.Range("K23").Value = 1.02688719799077
$ws.Range("L23").Value = 1.034010046548583
$ws.Range("M23").Value = 1.040945124589174
$ws.Range("N23").Value = 1.029424351324673
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023214725670281
$ws.Range("D24").Value = 1.025018503505079
$ws.Range("E24").Value = 1.032239317356519
$ws.Range("F24").Value = 1.039466910589367
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.029491518160384
$ws.Range("K24").Value = 1.02842733940359
$ws.Range("L24").Value = 1.035622554983063
$ws.Range("M24").Value = 1.042824918089508
$ws.Range("N24").Value = 1.030953514200499
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025605095095596
$ws.Range("D25").Value = 1.027131511033062
$ws.Range("E25").Value = 1.034434082541934
$ws.Range("F25").Value = 1.041971618222792
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.031259984173877
$ws.Range("K25").Value = 1.030211534456443
$ws.Range("L25").Value = 1.037490910712362
$ws.Range("M25").Value = 1.045004877348884
$ws.Range("N25").Value = 1.032724491638578